# Add a new "storageClassName" variable row to the sheet (chore: update storageclassname dinamico)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the last populated row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "storageClassName"
$ws.Cells.Item($newRow, 2).Value = "nfs-storage"
$ws.Cells.Item($newRow, 3).Value = "storage-nfs"

# Fill column D with the same shared concatenation formula used by the rows above
$ws.Cells.Item($newRow, 4).Formula = "=`$A`$1&`":`"&A$newRow&`"|`"&`$B`$1&`":`"&B$newRow&`"|`"&`$C`$1&`":`"&C$newRow"

# Match the selection left behind by the author's edit (one cell below/right of the new row)
$ws.Cells.Item($newRow + 1, 2).Select()
